$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the "value" figures that were originally recorded in units
# of 10,000 (万元). Re-scale them up by a factor of 10000 so the stored
# numbers represent plain units, consistent with the rest of the pipeline.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 * 10000
    }
}
